$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.580.96"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").Value = "3.318.14"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.41"
$ws.Range("E5").Value = "  -3.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.87"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "3.314.51"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("D13").Value = "3.894.35"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.12"
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("D16").Value = "66.646.48"
$ws.Range("E16").Value = "  -4.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000167"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").Value = "3.337.98"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "438.03"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.68"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.68"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.65"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.68"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "3.454.12"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.514"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.191"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  -6.66%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.96"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.87"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.32"
$ws.Range("E33").Value = "  -6.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.76"
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.18"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.28"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("D41").Value = "2.807.05"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.787"
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.46"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.23"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0676"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.17"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.22"
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  -8.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "319.43"
$ws.Range("E49").Value = "  -7.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0273"
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.984"
$ws.Range("E51").Value = "  -3.03%  "
